function Replace-InParagraph($doc, $paraIndex, $oldSub, $newSub) {
    # Replaces the first occurrence of $oldSub (plain logical text, can span
    # several runs / run-formatting boundaries) with $newSub inside paragraph
    # number $paraIndex, re-resolving the paragraph's Range each call so
    # offsets stay correct across several edits to the same paragraph.
    $p = $doc.Paragraphs.Item($paraIndex)
    $pStart = $p.Range.Start
    $text = $p.Range.Text
    $idx = $text.IndexOf($oldSub)
    if ($idx -lt 0) {
        throw "Replace-InParagraph: substring not found in paragraph $paraIndex : $oldSub"
    }
    $rngStart = $pStart + $idx
    $rngEnd = $rngStart + $oldSub.Length
    $rng = $doc.Range($rngStart, $rngEnd)
    $rng.Text = $newSub
}

function Collapse-RunsInParagraph($doc, $paraIndex, $span) {
    # $span's text does not change, but it is currently split across several
    # runs (e.g. a differently-styled "if" keyword run in the middle). Force
    # Word to rewrite/re-merge the run(s) covering $span into a single plain
    # run by round-tripping through a placeholder value (an in-place
    # assignment of identical text is treated as a no-op and would leave the
    # old run split as-is).
    $placeholder = "`u{E000}PLACEHOLDER`u{E000}"
    Replace-InParagraph $doc $paraIndex $span $placeholder
    Replace-InParagraph $doc $paraIndex $placeholder $span
}

$d = $word.ActiveDocument

# 1) Update the "Date" paragraph timestamp.
Replace-InParagraph $d 3 "June  28, 2021 (06:29:55 PM)" "June  28, 2021 (06:40:34 PM)"

# 2) "Let computer choose ..." -> "Starts by having computer choose ..."
Replace-InParagraph $d 33 `
    "Let computer choose a random number between 0 and 100 and store that number at a variable." `
    "Starts by having computer choose a random number between 0 and 100. Store that number at a variable."

# 3) "Asks the user to enter a numerical value, and stores ..." -> add "between 0 and 100"
Replace-InParagraph $d 34 `
    "Asks the user to enter a numerical value, and stores the user" `
    "Asks the user to enter a numerical value between 0 and 100 and stores the user"

# 4) "Add an if statement, displays on the screen" ("You guessed correctly" item):
#    the three runs ("Add an" / styled "if" / "statement, displays on the
#    screen") collapse into a single plain run with updated wording.
Replace-InParagraph $d 35 `
    "Add an if statement, displays on the screen" `
    "Add an if statement that displays on the screen"

# 5) The next three items ("Too high!", "Too low!", "multiple") already read
#    "Add an if statement that displays on the screen" but still have the
#    "if" keyword split into its own styled run - collapse each into one
#    plain run (text itself is unchanged here).
Collapse-RunsInParagraph $d 36 "Add an if statement that displays on the screen"
Collapse-RunsInParagraph $d 37 "Add an if statement that displays on the screen"
Collapse-RunsInParagraph $d 38 "Add an if statement that displays on the screen"

# 6) Fix typo: "is strictly than" -> "is strictly smaller than" (the "Too low!" item).
Replace-InParagraph $d 37 `
    "if the number entered by the user is strictly than the number selected by computer." `
    "if the number entered by the user is strictly smaller than the number selected by computer."

# 7) StringTok text for the last if-statement item: drop "of my favorite number".
Replace-InParagraph $d 38 `
    "You found a multiple of my favorite number!" `
    "You found a multiple !"

# 8) Closing paragraph: "... a multiple of your favorite number, ..." ->
#    "... a multiple of computer's number, ...".
Replace-InParagraph $d 42 `
    "a multiple of your favorite number, only one message is displayed." `
    "a multiple of computer’s number, only one message is displayed."
